$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 32000
$ws.Range("J13").Value = 32000
$ws.Range("L13").Value = 32000
$ws.Range("N13").Value = -32338
$ws.Range("H32").Value = 920
$ws.Range("I32").Value = 581.7778
$ws.Range("J32").Value = 1137.4286
$ws.Range("K32").Value = 581.7778
$ws.Range("L32").Value = 1137.4286
$ws.Range("M32").Value = -255.7778
$ws.Range("N32").Value = -1789.4286
$ws.Range("H53").Value = 454.04166
$ws.Range("I53").Value = 181.1875
$ws.Range("J53").Value = 999.75
$ws.Range("K53").Value = 181.1875
$ws.Range("L53").Value = 999.75
$ws.Range("M53").Value = 455.8125
$ws.Range("N53").Value = -2273.75
$ws.Range("H106").Value = 2090.1365
$ws.Range("I106").Value = 1666.7646
$ws.Range("J106").Value = 3529.6
$ws.Range("K106").Value = 1666.7646
$ws.Range("L106").Value = 3529.6
$ws.Range("M106").Value = -1035.7646
$ws.Range("N106").Value = -4791.6
$ws.Range("H107").Value = 1447.2222
$ws.Range("I107").Value = 1473.5294
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1473.5294
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 446.4706000000001
$ws.Range("N107").Value = -4840
$ws.Range("H132").Value = 21828676
$ws.Range("I132").Value = 25745012
$ws.Range("J132").Value = 9085.429
$ws.Range("K132").Value = 77235036
$ws.Range("L132").Value = 27256.287
$ws.Range("M132").Value = -77232506
$ws.Range("N132").Value = -32316.287
$ws.Range("H135").Value = 832.5714
$ws.Range("I135").Value = 832.5714
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 7493.1426
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -4958.1426
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4566.71
$ws.Range("I32").Value = 4670.132
$ws.Range("K32").Value = 4670.132
$ws.Range("M32").Value = -4383.132
$ws.Range("H97").Value = 982
$ws.Range("I97").Value = 829.1667
$ws.Range("K97").Value = 829.1667
$ws.Range("M97").Value = -333.1667
$ws.Range("H110").Value = 1642.75
$ws.Range("I110").Value = 1642.75
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1642.75
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 402.25
$ws.Range("N110").ClearContents()
$ws.Range("H137").Value = 41780
$ws.Range("J137").Value = 41780
$ws.Range("L137").Value = 41780
$ws.Range("N137").Value = -51980

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 7517056
$ws.Range("I7").Value = 11676668
$ws.Range("J7").Value = 5734365.5
$ws.Range("K7").Value = 11676668
$ws.Range("L7").Value = 5734365.5
$ws.Range("M7").Value = -11676555
$ws.Range("N7").Value = -5734591.5
$ws.Range("H20").Value = 17298.385
$ws.Range("I20").Value = 9276.2
$ws.Range("J20").Value = 22312.25
$ws.Range("K20").Value = 9276.2
$ws.Range("L20").Value = 22312.25
$ws.Range("M20").Value = -9029.2
$ws.Range("N20").Value = -22806.25
$ws.Range("H107").Value = 2665.3333
$ws.Range("I107").Value = 2798.4
$ws.Range("K107").Value = 2798.4
$ws.Range("M107").Value = -878.4000000000001
$ws.Range("H137").Value = 49827.5
$ws.Range("J137").Value = 49827.5
$ws.Range("L137").Value = 49827.5
$ws.Range("N137").Value = -60027.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 235515.27
$ws.Range("I31").Value = 410472.34
$ws.Range("J31").Value = 4571.92
$ws.Range("K31").Value = 410472.34
$ws.Range("L31").Value = 4571.92
$ws.Range("M31").Value = -410177.34
$ws.Range("N31").Value = -5161.92
$ws.Range("H34").Value = 235515.27
$ws.Range("I34").Value = 410472.34
$ws.Range("J34").Value = 4571.92
$ws.Range("K34").Value = 410472.34
$ws.Range("L34").Value = 4571.92
$ws.Range("M34").Value = -410270.34
$ws.Range("N34").Value = -4975.92
$ws.Range("H105").Value = 2884.36
$ws.Range("I105").Value = 3842
$ws.Range("J105").Value = 2644.95
$ws.Range("K105").Value = 3842
$ws.Range("L105").Value = 2644.95
$ws.Range("M105").Value = -2095
$ws.Range("N105").Value = -6138.95
$ws.Range("H134").Value = 6116.773
$ws.Range("I134").Value = 5903.45
$ws.Range("J134").Value = 8250
$ws.Range("K134").Value = 17710.35
$ws.Range("L134").Value = 24750
$ws.Range("M134").Value = -15175.35
$ws.Range("N134").Value = -29820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 312.3
$ws.Range("I50").Value = 142.42857
$ws.Range("J50").Value = 708.6667
$ws.Range("K50").Value = 427.28571
$ws.Range("L50").Value = 2126.0001
$ws.Range("M50").Value = 53.71429000000001
$ws.Range("N50").Value = -3088.0001
$ws.Range("H53").Value = 312.3
$ws.Range("I53").Value = 142.42857
$ws.Range("J53").Value = 708.6667
$ws.Range("K53").Value = 427.28571
$ws.Range("L53").Value = 2126.0001
$ws.Range("M53").Value = 53.71429000000001
$ws.Range("N53").Value = -3088.0001
$ws.Range("H68").Value = 2864.7188
$ws.Range("I68").Value = 1330
$ws.Range("J68").Value = 3294.44
$ws.Range("K68").Value = 3990
$ws.Range("L68").Value = 9883.32
$ws.Range("M68").Value = -3179
$ws.Range("N68").Value = -11505.32
$ws.Range("H71").Value = 2864.7188
$ws.Range("I71").Value = 1330
$ws.Range("J71").Value = 3294.44
$ws.Range("K71").Value = 11970
$ws.Range("L71").Value = 29649.96
$ws.Range("M71").Value = -7914
$ws.Range("N71").Value = -37761.96
$ws.Range("H113").Value = 484.25
$ws.Range("I113").Value = 481
$ws.Range("J113").Value = 490
$ws.Range("K113").Value = 1443
$ws.Range("L113").Value = 1470
$ws.Range("M113").Value = 727
$ws.Range("N113").Value = -5810
$ws.Range("H131").Value = 760.87
$ws.Range("J131").Value = 806.3483
$ws.Range("L131").Value = 2419.0449
$ws.Range("N131").Value = -12499.0449

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 24004.154
$ws.Range("J46").Value = 24754.916
$ws.Range("L46").Value = 24754.916
$ws.Range("N46").Value = -25066.916
$ws.Range("H128").Value = 39221.125
$ws.Range("J128").Value = 41865.715
$ws.Range("L128").Value = 41865.715
$ws.Range("N128").Value = -51825.715
$ws.Range("H137").Value = 39113.332
$ws.Range("J137").Value = 48670
$ws.Range("L137").Value = 48670
$ws.Range("N137").Value = -58870

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 612.63635
$ws.Range("I16").Value = 526.6667
$ws.Range("J16").Value = 999.5
$ws.Range("K16").Value = 526.6667
$ws.Range("L16").Value = 999.5
$ws.Range("M16").Value = -356.6667
$ws.Range("N16").Value = -1339.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2562.3794
$ws.Range("I132").Value = 831.61536
$ws.Range("J132").Value = 3968.625
$ws.Range("K132").Value = 2494.84608
$ws.Range("L132").Value = 11905.875
$ws.Range("M132").Value = 35.15391999999974
$ws.Range("N132").Value = -16965.875
$ws.Range("H136").Value = 3818.1904
$ws.Range("I136").Value = 775.41174
$ws.Range("J136").Value = 16750
$ws.Range("K136").Value = 2326.23522
$ws.Range("L136").Value = 50250
$ws.Range("M136").Value = 223.76478
$ws.Range("N136").Value = -55350
